# Update Active_Outages.xlsx - 6/18/2025, 5:00:06 PM
$wb = $excel.ActiveWorkbook

# --- Sheet R1: elapsed duration refresh (rows 2-4), and row 5 replaced with new pending site ---
$ws1 = $wb.Worksheets.Item("R1")
$ws1.Range("G2").Value = "3930:14:12"
$ws1.Range("G3").Value = "69:46:50"
$ws1.Range("G4").Value = "92:46:50"
$ws1.Range("B5").Value = "R4"
$ws1.Range("D5").Value = "JED0125"
$ws1.Range("I5").Value = "Generator-SG"
$ws1.Range("J5").Value = "Good+In progress"

# --- Sheet R2: elapsed duration refresh (rows 2-4) ---
$ws2 = $wb.Worksheets.Item("R2")
$ws2.Range("G2").Value = "12111:37:53"
$ws2.Range("G3").Value = "3241:21:22"
$ws2.Range("G4").Value = "479:32:56"

# --- Sheet R4: elapsed duration refresh (rows 2-5) ---
$ws4 = $wb.Worksheets.Item("R4")
$ws4.Range("G2").Value = "2957:27:42"
$ws4.Range("G3").Value = "184:39:57"
$ws4.Range("G4").Value = "72:52:22"
$ws4.Range("G5").Value = "70:29:55"

# --- Sheet R5: elapsed duration refresh (row 2) ---
$ws5 = $wb.Worksheets.Item("R5")
$ws5.Range("G2").Value = "431:26:41"

# --- Sheet R6: elapsed duration refresh (row 2) ---
$ws6 = $wb.Worksheets.Item("R6")
$ws6.Range("G2").Value = "71:58:59"
